$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$meta.Range("B3").Value = "1.8.2"

# Date: 2023-06-27T22:42:19-04:00 -> 2023-09-01T14:45:29-04:00
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 1 is the root "Extension" element. Its Invariants cell (column AJ) was blank
# and now carries the ele-1 / ext-1 constraint text (same text already used on the
# Extension.extension row).
$invariantText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Range("AJ1").Value = $invariantText
